$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = "2025/12/05 01:00"
$ws.Range("B56").Value = "-"
$ws.Range("C56").Value = "-"
$ws.Range("D56").Value = "-"
$ws.Range("E56").Value = "-"
$ws.Range("F56").Value = "-"
$ws.Range("G56").Value = "-"
